# update methanol reactor data
# Workbook: Model_Data_Base.xlsx - "Units" table (Table1, A1:AJ6)
# Columns touched: C=Input2, J=Cap_Output1_existing, L=Cap_Output2_existing,
#                   S=start_up_Output1, U=shut_down_Output1, W=Relation_In_In,
#                   Y=Relation_Out_Out, AA=unit_on_cost, AB=fom_cost

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# Row 5 - Destilation_Tower: give it a second input (Steam) and a Relation_In_In value
$ws.Range("C5").Value = "Steam"
$ws.Range("W5").Value = 11.03

# Row 6 - Methanol_Reactor: revised capacities / ramp / relation figures
$ws.Range("J6").Value = 52          # Cap_Output1_existing: 100 -> 52
$ws.Range("L6").ClearContents()     # Cap_Output2_existing: 100 -> (blank)
$ws.Range("S6").Value = 0.5         # start_up_Output1: (blank) -> 0.5
$ws.Range("U6").Value = 0.5         # shut_down_Output1: (blank) -> 0.5
$ws.Range("W6").Value = 4.57        # Relation_In_In: 1 -> 4.57
$ws.Range("Y6").Value = 4.32        # Relation_Out_Out: 4 -> 4.32
$ws.Range("AA6").ClearContents()    # unit_on_cost: 0.0000001 -> (blank)
$ws.Range("AB6").Value = 4.45       # fom_cost: (blank) -> 4.45

# Update the sheet's saved selection/scroll position
$ws.Range("L9").Select()
